$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the header style (bold, border, centered) from F1 to G1:H1
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update existing B2 and D2 values with new precision
$ws.Range("B2").Value = 0.07514644587374582
$ws.Range("D2").Value = 0.2119198634755614

# Add new data cells
$ws.Range("G2").Value = 0.1258822953001072
$ws.Range("H2").Value = 0.988
